$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.790.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.701.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'317.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.4090"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.006"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'52.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08930"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.731"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.92%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'24.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'8.184"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.00001334"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.713.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'99.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.07156"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'20.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.265"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.51%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'14.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'24.791.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.108"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.338"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'23.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.285"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +23.27%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'165.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'139.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'5.198"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'8.177"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.32%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09168"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.60%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.081"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.03061"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.33%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.2822"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'14.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.09313"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.7835"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'1.475"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'16.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.647"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.7265"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'4.242"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.365"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'141.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'93.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.25%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.08055"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.69%  "
$ws.Range("E51").Style = "Normal"
